$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 21999.5  # H3: 22328.5 -> 21999.5
$ws.Cells.Item(3, 10).Value = 21999.5  # J3: 22328.5 -> 21999.5
$ws.Cells.Item(3, 12).Value = 21999.5  # L3: 22328.5 -> 21999.5
$ws.Cells.Item(3, 14).Value = -22227.5  # N3: -22556.5 -> -22227.5
$ws.Cells.Item(38, 8).Value = 894  # H38: 1250674.5 -> 894
$ws.Cells.Item(38, 9).Value = 692  # I38: 2000179.2 -> 692
$ws.Cells.Item(38, 11).Value = 2076  # K38: 6000537.6 -> 2076
$ws.Cells.Item(38, 13).Value = -1704  # M38: -6000165.6 -> -1704
$ws.Cells.Item(40, 8).Value = 1616.3889  # H40: 1589.1052 -> 1616.3889
$ws.Cells.Item(40, 10).Value = 2249.75  # J40: 2121.7778 -> 2249.75
$ws.Cells.Item(40, 12).Value = 2249.75  # L40: 2121.7778 -> 2249.75
$ws.Cells.Item(40, 14).Value = -2599.75  # N40: -2471.7778 -> -2599.75
$ws.Cells.Item(43, 8).Value = 8004.5  # H43: 8003 -> 8004.5
$ws.Cells.Item(43, 9).Value = 8001  # I43: 8000.5 -> 8001
$ws.Cells.Item(43, 11).Value = 8001  # K43: 8000.5 -> 8001
$ws.Cells.Item(43, 13).Value = -7932  # M43: -7931.5 -> -7932
$ws.Cells.Item(51, 8).Value = 7211  # H51: 6719.8 -> 7211
$ws.Cells.Item(51, 10).Value = 12999.75  # J51: 10859.6 -> 12999.75
$ws.Cells.Item(51, 12).Value = 12999.75  # L51: 10859.6 -> 12999.75
$ws.Cells.Item(51, 14).Value = -13967.75  # N51: -11827.6 -> -13967.75
$ws.Cells.Item(58, 8).Value = 527.5  # H58: 2013.75 -> 527.5
$ws.Cells.Item(58, 10).Value = 0  # J58: 3500 -> 0
$ws.Cells.Item(58, 12).Value = 0  # L58: 10500 -> 0
$ws.Cells.Item(58, 14).ClearContents()  # N58: -10800 -> (removed)
$ws.Cells.Item(70, 8).Value = 2000  # H70: 2425.25 -> 2000
$ws.Cells.Item(70, 9).Value = 0  # I70: 2202 -> 0
$ws.Cells.Item(70, 10).Value = 2000  # J70: 2499.6667 -> 2000
$ws.Cells.Item(70, 11).Value = 0  # K70: 6606 -> 0
$ws.Cells.Item(70, 12).Value = 6000  # L70: 7499.000100000001 -> 6000
$ws.Cells.Item(70, 13).ClearContents()  # M70: -6336 -> (removed)
$ws.Cells.Item(70, 14).Value = -6540  # N70: -8039.000100000001 -> -6540
$ws.Cells.Item(73, 8).Value = 2000  # H73: 2425.25 -> 2000
$ws.Cells.Item(73, 9).Value = 0  # I73: 2202 -> 0
$ws.Cells.Item(73, 10).Value = 2000  # J73: 2499.6667 -> 2000
$ws.Cells.Item(73, 11).Value = 0  # K73: 6606 -> 0
$ws.Cells.Item(73, 12).Value = 6000  # L73: 7499.000100000001 -> 6000
$ws.Cells.Item(73, 13).ClearContents()  # M73: -5670 -> (removed)
$ws.Cells.Item(73, 14).Value = -7872  # N73: -9371.000100000001 -> -7872
$ws.Cells.Item(74, 8).Value = 4950  # H74: 3500 -> 4950
$ws.Cells.Item(74, 9).Value = 4950  # I74: 3500 -> 4950
$ws.Cells.Item(74, 11).Value = 4950  # K74: 3500 -> 4950
$ws.Cells.Item(74, 13).Value = -4014  # M74: -2564 -> -4014
$ws.Cells.Item(77, 8).Value = 4950  # H77: 3500 -> 4950
$ws.Cells.Item(77, 9).Value = 4950  # I77: 3500 -> 4950
$ws.Cells.Item(77, 11).Value = 24750  # K77: 17500 -> 24750
$ws.Cells.Item(77, 13).Value = -20070  # M77: -12820 -> -20070
$ws.Cells.Item(86, 8).Value = 3049.875  # H86: 3342.8572 -> 3049.875
$ws.Cells.Item(86, 10).Value = 2799.6667  # J86: 3700 -> 2799.6667
$ws.Cells.Item(86, 12).Value = 2799.6667  # L86: 3700 -> 2799.6667
$ws.Cells.Item(86, 14).Value = -5045.6667  # N86: -5946 -> -5045.6667
$ws.Cells.Item(87, 8).Value = 0  # H87: 89999 -> 0
$ws.Cells.Item(87, 10).Value = 0  # J87: 89999 -> 0
$ws.Cells.Item(87, 12).Value = 0  # L87: 89999 -> 0
$ws.Cells.Item(87, 14).ClearContents()  # N87: -92495 -> (removed)
$ws.Cells.Item(89, 8).Value = 3049.875  # H89: 3342.8572 -> 3049.875
$ws.Cells.Item(89, 10).Value = 2799.6667  # J89: 3700 -> 2799.6667
$ws.Cells.Item(89, 12).Value = 13998.3335  # L89: 18500 -> 13998.3335
$ws.Cells.Item(89, 14).Value = -25230.3335  # N89: -29732 -> -25230.3335
$ws.Cells.Item(90, 8).Value = 0  # H90: 89999 -> 0
$ws.Cells.Item(90, 10).Value = 0  # J90: 89999 -> 0
$ws.Cells.Item(90, 12).Value = 0  # L90: 269997 -> 0
$ws.Cells.Item(90, 14).ClearContents()  # N90: -282477 -> (removed)
$ws.Cells.Item(100, 8).Value = 2711  # H100: 2399.4546 -> 2711
$ws.Cells.Item(100, 9).Value = 2666.8333  # I100: 2428.4285 -> 2666.8333
$ws.Cells.Item(100, 10).Value = 2799.3333  # J100: 2348.75 -> 2799.3333
$ws.Cells.Item(100, 11).Value = 2666.8333  # K100: 2428.4285 -> 2666.8333
$ws.Cells.Item(100, 12).Value = 2799.3333  # L100: 2348.75 -> 2799.3333
$ws.Cells.Item(100, 13).Value = -2125.8333  # M100: -1887.4285 -> -2125.8333
$ws.Cells.Item(100, 14).Value = -3881.3333  # N100: -3430.75 -> -3881.3333
$ws.Cells.Item(102, 8).Value = 21999.5  # H102: 22328.5 -> 21999.5
$ws.Cells.Item(102, 10).Value = 21999.5  # J102: 22328.5 -> 21999.5
$ws.Cells.Item(102, 12).Value = 21999.5  # L102: 22328.5 -> 21999.5
$ws.Cells.Item(102, 14).Value = -28489.5  # N102: -28818.5 -> -28489.5
$ws.Cells.Item(106, 8).Value = 8850.75  # H106: 8931.477000000001 -> 8850.75
$ws.Cells.Item(106, 9).Value = 6278.778  # I106: 6503.3687 -> 6278.778
$ws.Cells.Item(106, 11).Value = 6278.778  # K106: 6503.3687 -> 6278.778
$ws.Cells.Item(106, 13).Value = -5647.778  # M106: -5872.3687 -> -5647.778
$ws.Cells.Item(107, 8).Value = 1567.1818  # H107: 1593 -> 1567.1818
$ws.Cells.Item(107, 9).Value = 1554  # I107: 1595.5714 -> 1554
$ws.Cells.Item(107, 10).Value = 1590.25  # J107: 1587 -> 1590.25
$ws.Cells.Item(107, 11).Value = 1554  # K107: 1595.5714 -> 1554
$ws.Cells.Item(107, 12).Value = 1590.25  # L107: 1587 -> 1590.25
$ws.Cells.Item(107, 13).Value = 366  # M107: 324.4286 -> 366
$ws.Cells.Item(107, 14).Value = -5430.25  # N107: -5427 -> -5430.25
$ws.Cells.Item(111, 8).Value = 2721.5  # H111: 2999.5 -> 2721.5
$ws.Cells.Item(111, 9).Value = 5000  # I111: 2999.5 -> 5000
$ws.Cells.Item(111, 10).Value = 443  # J111: 0 -> 443
$ws.Cells.Item(111, 11).Value = 15000  # K111: 8998.5 -> 15000
$ws.Cells.Item(111, 12).Value = 1329  # L111: 0 -> 1329
$ws.Cells.Item(111, 13).Value = -11933  # M111: -5931.5 -> -11933
$ws.Cells.Item(111, 14).Value = -7463  # N111: None -> -7463
$ws.Cells.Item(113, 8).Value = 8628.143  # H113: 8642.429 -> 8628.143
$ws.Cells.Item(113, 10).Value = 9599.4  # J113: 9619.4 -> 9599.4
$ws.Cells.Item(113, 12).Value = 9599.4  # L113: 9619.4 -> 9599.4
$ws.Cells.Item(113, 14).Value = -16107.4  # N113: -16127.4 -> -16107.4
$ws.Cells.Item(116, 8).Value = 4938.857  # H116: 4945.5 -> 4938.857
$ws.Cells.Item(116, 9).Value = 4945.5  # I116: 4952.143 -> 4945.5
$ws.Cells.Item(116, 11).Value = 4945.5  # K116: 4952.143 -> 4945.5
$ws.Cells.Item(116, 13).Value = -1503.5  # M116: -1510.143 -> -1503.5
$ws.Cells.Item(118, 8).Value = 340  # H118: 323 -> 340
$ws.Cells.Item(118, 9).Value = 340  # I118: 323 -> 340
$ws.Cells.Item(118, 11).Value = 1020  # K118: 969 -> 1020
$ws.Cells.Item(118, 13).Value = 637  # M118: 688 -> 637
$ws.Cells.Item(125, 8).Value = 0  # H125: 583.3333 -> 0
$ws.Cells.Item(125, 9).Value = 0  # I125: 750 -> 0
$ws.Cells.Item(125, 10).Value = 0  # J125: 500 -> 0
$ws.Cells.Item(125, 11).Value = 0  # K125: 6750 -> 0
$ws.Cells.Item(125, 12).Value = 0  # L125: 4500 -> 0
$ws.Cells.Item(125, 13).ClearContents()  # M125: -4290 -> (removed)
$ws.Cells.Item(125, 14).ClearContents()  # N125: -9420 -> (removed)
$ws.Cells.Item(131, 8).Value = 713.1111  # H131: 719.7778 -> 713.1111
$ws.Cells.Item(131, 9).Value = 755.1429000000001  # I131: 763.7143 -> 755.1429000000001
$ws.Cells.Item(131, 11).Value = 2265.4287  # K131: 2291.1429 -> 2265.4287
$ws.Cells.Item(131, 13).Value = 2774.5713  # M131: 2748.8571 -> 2774.5713
$ws.Cells.Item(132, 8).Value = 2495.611  # H132: 2554.2942 -> 2495.611
$ws.Cells.Item(132, 9).Value = 2495.611  # I132: 2554.2942 -> 2495.611
$ws.Cells.Item(132, 11).Value = 7486.833  # K132: 7662.882599999999 -> 7486.833
$ws.Cells.Item(132, 13).Value = -4956.833  # M132: -5132.882599999999 -> -4956.833
$ws.Cells.Item(137, 8).Value = 4999  # H137: 2782.3333 -> 4999
$ws.Cells.Item(137, 9).Value = 0  # I137: 2339 -> 0
$ws.Cells.Item(137, 11).Value = 0  # K137: 7017 -> 0
$ws.Cells.Item(137, 13).ClearContents()  # M137: -4467 -> (removed)
$ws.Cells.Item(138, 8).Value = 2211.6667  # H138: 1843.3334 -> 2211.6667
$ws.Cells.Item(138, 9).Value = 818.125  # I138: 941.5714 -> 818.125
$ws.Cells.Item(138, 10).Value = 4998.75  # J138: 4999.5 -> 4998.75
$ws.Cells.Item(138, 11).Value = 2454.375  # K138: 2824.7142 -> 2454.375
$ws.Cells.Item(138, 12).Value = 14996.25  # L138: 14998.5 -> 14996.25
$ws.Cells.Item(138, 13).Value = 2685.625  # M138: 2315.2858 -> 2685.625
$ws.Cells.Item(138, 14).Value = -25276.25  # N138: -25278.5 -> -25276.25
$ws.Cells.Item(141, 8).Value = 13360.25  # H141: 15112.714 -> 13360.25
$ws.Cells.Item(141, 9).Value = 14697.5  # I141: 16923.75 -> 14697.5
$ws.Cells.Item(141, 10).Value = 12023  # J141: 12698 -> 12023
$ws.Cells.Item(141, 11).Value = 44092.5  # K141: 50771.25 -> 44092.5
$ws.Cells.Item(141, 12).Value = 36069  # L141: 38094 -> 36069
$ws.Cells.Item(141, 13).Value = -38912.5  # M141: -45591.25 -> -38912.5
$ws.Cells.Item(141, 14).Value = -46429  # N141: -48454 -> -46429

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 1266  # H4: 1300 -> 1266
$ws.Cells.Item(4, 9).Value = 1266  # I4: 1300 -> 1266
$ws.Cells.Item(4, 11).Value = 1266  # K4: 1300 -> 1266
$ws.Cells.Item(4, 13).Value = -1150  # M4: -1184 -> -1150
$ws.Cells.Item(16, 8).Value = 39999  # H16: 20499.5 -> 39999
$ws.Cells.Item(16, 9).Value = 0  # I16: 1000 -> 0
$ws.Cells.Item(16, 11).Value = 0  # K16: 1000 -> 0
$ws.Cells.Item(16, 13).ClearContents()  # M16: -713 -> (removed)
$ws.Cells.Item(32, 8).Value = 21900.092  # H32: 20324.25 -> 21900.092
$ws.Cells.Item(32, 9).Value = 23790.2  # I32: 21899.273 -> 23790.2
$ws.Cells.Item(32, 11).Value = 23790.2  # K32: 21899.273 -> 23790.2
$ws.Cells.Item(32, 13).Value = -23503.2  # M32: -21612.273 -> -23503.2
$ws.Cells.Item(34, 8).Value = 0  # H34: 50000000 -> 0
$ws.Cells.Item(34, 9).Value = 0  # I34: 50000000 -> 0
$ws.Cells.Item(34, 11).Value = 0  # K34: 50000000 -> 0
$ws.Cells.Item(34, 13).ClearContents()  # M34: -49999729 -> (removed)
$ws.Cells.Item(40, 8).Value = 0  # H40: 15800 -> 0
$ws.Cells.Item(40, 10).Value = 0  # J40: 15800 -> 0
$ws.Cells.Item(40, 12).Value = 0  # L40: 15800 -> 0
$ws.Cells.Item(40, 14).ClearContents()  # N40: -16152 -> (removed)
$ws.Cells.Item(45, 8).Value = 1438.125  # H45: 1649.625 -> 1438.125
$ws.Cells.Item(45, 9).Value = 1438.125  # I45: 1613.8572 -> 1438.125
$ws.Cells.Item(45, 10).Value = 0  # J45: 1900 -> 0
$ws.Cells.Item(45, 11).Value = 1438.125  # K45: 1613.8572 -> 1438.125
$ws.Cells.Item(45, 12).Value = 0  # L45: 1900 -> 0
$ws.Cells.Item(45, 13).Value = -1061.125  # M45: -1236.8572 -> -1061.125
$ws.Cells.Item(45, 14).ClearContents()  # N45: -2654 -> (removed)
$ws.Cells.Item(61, 8).Value = 1441.5834  # H61: 1499.909 -> 1441.5834
$ws.Cells.Item(61, 10).Value = 800  # J61: 0 -> 800
$ws.Cells.Item(61, 12).Value = 800  # L61: 0 -> 800
$ws.Cells.Item(61, 14).Value = -1224  # N61: None -> -1224
$ws.Cells.Item(74, 8).Value = 1768.1  # H74: 1785 -> 1768.1
$ws.Cells.Item(74, 9).Value = 1768.5  # I74: 1792.4 -> 1768.5
$ws.Cells.Item(74, 11).Value = 1768.5  # K74: 1792.4 -> 1768.5
$ws.Cells.Item(74, 13).Value = -894.5  # M74: -918.4000000000001 -> -894.5
$ws.Cells.Item(77, 8).Value = 1768.1  # H77: 1785 -> 1768.1
$ws.Cells.Item(77, 9).Value = 1768.5  # I77: 1792.4 -> 1768.5
$ws.Cells.Item(77, 11).Value = 8842.5  # K77: 8962 -> 8842.5
$ws.Cells.Item(77, 13).Value = -4474.5  # M77: -4594 -> -4474.5
$ws.Cells.Item(88, 8).Value = 1185.5625  # H88: 1413.4166 -> 1185.5625
$ws.Cells.Item(88, 9).Value = 810.375  # I88: 945.5 -> 810.375
$ws.Cells.Item(88, 10).Value = 1560.75  # J88: 1881.3334 -> 1560.75
$ws.Cells.Item(88, 11).Value = 810.375  # K88: 945.5 -> 810.375
$ws.Cells.Item(88, 12).Value = 1560.75  # L88: 1881.3334 -> 1560.75
$ws.Cells.Item(88, 13).Value = -404.375  # M88: -539.5 -> -404.375
$ws.Cells.Item(88, 14).Value = -2372.75  # N88: -2693.3334 -> -2372.75
$ws.Cells.Item(91, 8).Value = 1185.5625  # H91: 1413.4166 -> 1185.5625
$ws.Cells.Item(91, 9).Value = 810.375  # I91: 945.5 -> 810.375
$ws.Cells.Item(91, 10).Value = 1560.75  # J91: 1881.3334 -> 1560.75
$ws.Cells.Item(91, 11).Value = 810.375  # K91: 945.5 -> 810.375
$ws.Cells.Item(91, 12).Value = 1560.75  # L91: 1881.3334 -> 1560.75
$ws.Cells.Item(91, 13).Value = 593.625  # M91: 458.5 -> 593.625
$ws.Cells.Item(91, 14).Value = -4368.75  # N91: -4689.3334 -> -4368.75
$ws.Cells.Item(102, 8).Value = 7520455.5  # H102: 1818.7646 -> 7520455.5
$ws.Cells.Item(102, 9).Value = 8930173  # I102: 1787.2142 -> 8930173
$ws.Cells.Item(102, 10).Value = 1962.6666  # J102: 1966 -> 1962.6666
$ws.Cells.Item(102, 11).Value = 8930173  # K102: 1787.2142 -> 8930173
$ws.Cells.Item(102, 12).Value = 1962.6666  # L102: 1966 -> 1962.6666
$ws.Cells.Item(102, 13).Value = -8928551  # M102: -165.2141999999999 -> -8928551
$ws.Cells.Item(102, 14).Value = -5206.6666  # N102: -5210 -> -5206.6666
$ws.Cells.Item(125, 8).Value = 100238.336  # H125: 101000 -> 100238.336
$ws.Cells.Item(125, 10).Value = 100238.336  # J125: 101000 -> 100238.336
$ws.Cells.Item(125, 12).Value = 100238.336  # L125: 101000 -> 100238.336
$ws.Cells.Item(125, 14).Value = -110078.336  # N125: -110840 -> -110078.336
$ws.Cells.Item(132, 8).Value = 4351.273  # H132: 5177 -> 4351.273
$ws.Cells.Item(132, 9).Value = 1666.6666  # I132: 1500 -> 1666.6666
$ws.Cells.Item(132, 10).Value = 5358  # J132: 6157.533 -> 5358
$ws.Cells.Item(132, 11).Value = 4999.9998  # K132: 4500 -> 4999.9998
$ws.Cells.Item(132, 12).Value = 16074  # L132: 18472.599 -> 16074
$ws.Cells.Item(132, 13).Value = -2469.9998  # M132: -1970 -> -2469.9998
$ws.Cells.Item(132, 14).Value = -21134  # N132: -23532.599 -> -21134
$ws.Cells.Item(136, 8).Value = 1441.5834  # H136: 1499.909 -> 1441.5834
$ws.Cells.Item(136, 10).Value = 800  # J136: 0 -> 800
$ws.Cells.Item(136, 12).Value = 2400  # L136: 0 -> 2400
$ws.Cells.Item(136, 14).Value = -7500  # N136: None -> -7500

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(19, 8).Value = 0  # H19: 1200 -> 0
$ws.Cells.Item(19, 10).Value = 0  # J19: 1200 -> 0
$ws.Cells.Item(19, 12).Value = 0  # L19: 1200 -> 0
$ws.Cells.Item(19, 14).ClearContents()  # N19: -1546 -> (removed)
$ws.Cells.Item(86, 8).Value = 2901.25  # H86: 2913.6875 -> 2901.25
$ws.Cells.Item(86, 9).Value = 3041.5334  # I86: 2913.6875 -> 3041.5334
$ws.Cells.Item(86, 10).Value = 797  # J86: 0 -> 797
$ws.Cells.Item(86, 11).Value = 3041.5334  # K86: 2913.6875 -> 3041.5334
$ws.Cells.Item(86, 12).Value = 797  # L86: 0 -> 797
$ws.Cells.Item(86, 13).Value = -1918.5334  # M86: -1790.6875 -> -1918.5334
$ws.Cells.Item(86, 14).Value = -3043  # N86: None -> -3043
$ws.Cells.Item(89, 8).Value = 2901.25  # H89: 2913.6875 -> 2901.25
$ws.Cells.Item(89, 9).Value = 3041.5334  # I89: 2913.6875 -> 3041.5334
$ws.Cells.Item(89, 10).Value = 797  # J89: 0 -> 797
$ws.Cells.Item(89, 11).Value = 15207.667  # K89: 14568.4375 -> 15207.667
$ws.Cells.Item(89, 12).Value = 3985  # L89: 0 -> 3985
$ws.Cells.Item(89, 13).Value = -9591.666999999999  # M89: -8952.4375 -> -9591.666999999999
$ws.Cells.Item(89, 14).Value = -15217  # N89: None -> -15217
$ws.Cells.Item(94, 8).Value = 71429000  # H94: 454.7143 -> 71429000
$ws.Cells.Item(94, 9).Value = 250000260  # I94: 500 -> 250000260
$ws.Cells.Item(94, 10).Value = 498.8  # J94: 447.16666 -> 498.8
$ws.Cells.Item(94, 11).Value = 250000260  # K94: 500 -> 250000260
$ws.Cells.Item(94, 12).Value = 498.8  # L94: 447.16666 -> 498.8
$ws.Cells.Item(94, 13).Value = -249999809  # M94: -49 -> -249999809
$ws.Cells.Item(94, 14).Value = -1400.8  # N94: -1349.16666 -> -1400.8
$ws.Cells.Item(105, 8).Value = 66670776  # H105: 4402.857 -> 66670776
$ws.Cells.Item(105, 10).Value = 76927360  # J105: 4644.1665 -> 76927360
$ws.Cells.Item(105, 12).Value = 76927360  # L105: 4644.1665 -> 76927360
$ws.Cells.Item(105, 14).Value = -76930854  # N105: -8138.1665 -> -76930854
$ws.Cells.Item(107, 8).Value = 2828.7144  # H107: 2467 -> 2828.7144
$ws.Cells.Item(107, 9).Value = 2466.8333  # I107: 1960.4 -> 2466.8333
$ws.Cells.Item(107, 11).Value = 2466.8333  # K107: 1960.4 -> 2466.8333
$ws.Cells.Item(107, 13).Value = -546.8332999999998  # M107: -40.40000000000009 -> -546.8332999999998
$ws.Cells.Item(134, 8).Value = 1487.5  # H134: 1499.75 -> 1487.5
$ws.Cells.Item(134, 9).Value = 1487.5  # I134: 1499.75 -> 1487.5
$ws.Cells.Item(134, 11).Value = 4462.5  # K134: 4499.25 -> 4462.5
$ws.Cells.Item(134, 13).Value = -1927.5  # M134: -1964.25 -> -1927.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 2999  # H4: 3999 -> 2999
$ws.Cells.Item(4, 10).Value = 0  # J4: 6999 -> 0
$ws.Cells.Item(4, 12).Value = 0  # L4: 6999 -> 0
$ws.Cells.Item(4, 14).ClearContents()  # N4: -7223 -> (removed)
$ws.Cells.Item(22, 8).Value = 14933.889  # H22: 14933.777 -> 14933.889
$ws.Cells.Item(22, 9).Value = 3100  # I22: 3950 -> 3100
$ws.Cells.Item(22, 10).Value = 20850.834  # J22: 18072 -> 20850.834
$ws.Cells.Item(22, 11).Value = 3100  # K22: 3950 -> 3100
$ws.Cells.Item(22, 12).Value = 20850.834  # L22: 18072 -> 20850.834
$ws.Cells.Item(22, 13).Value = -2750  # M22: -3600 -> -2750
$ws.Cells.Item(22, 14).Value = -21550.834  # N22: -18772 -> -21550.834
$ws.Cells.Item(31, 8).Value = 1916.875  # H31: 1954.25 -> 1916.875
$ws.Cells.Item(31, 9).Value = 1833.7142  # I31: 1876.4286 -> 1833.7142
$ws.Cells.Item(31, 11).Value = 1833.7142  # K31: 1876.4286 -> 1833.7142
$ws.Cells.Item(31, 13).Value = -1538.7142  # M31: -1581.4286 -> -1538.7142
$ws.Cells.Item(34, 8).Value = 1916.875  # H34: 1954.25 -> 1916.875
$ws.Cells.Item(34, 9).Value = 1833.7142  # I34: 1876.4286 -> 1833.7142
$ws.Cells.Item(34, 11).Value = 1833.7142  # K34: 1876.4286 -> 1833.7142
$ws.Cells.Item(34, 13).Value = -1631.7142  # M34: -1674.4286 -> -1631.7142
$ws.Cells.Item(132, 8).Value = 2472  # H132: 2531.125 -> 2472
$ws.Cells.Item(132, 10).Value = 1999  # J132: 0 -> 1999
$ws.Cells.Item(132, 12).Value = 5997  # L132: 0 -> 5997
$ws.Cells.Item(132, 14).Value = -11057  # N132: None -> -11057
$ws.Cells.Item(134, 8).Value = 2947.4443  # H134: 3103.3333 -> 2947.4443
$ws.Cells.Item(134, 9).Value = 2988.4707  # I134: 3164.2856 -> 2988.4707
$ws.Cells.Item(134, 11).Value = 8965.4121  # K134: 9492.856800000001 -> 8965.4121
$ws.Cells.Item(134, 13).Value = -6430.4121  # M134: -6957.856800000001 -> -6430.4121

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 1594.3334  # H11: 2063.3333 -> 1594.3334
$ws.Cells.Item(11, 9).Value = 187.33333  # I11: 188.5 -> 187.33333
$ws.Cells.Item(11, 10).Value = 3001.3333  # J11: 3000.75 -> 3001.3333
$ws.Cells.Item(11, 11).Value = 561.99999  # K11: 565.5 -> 561.99999
$ws.Cells.Item(11, 12).Value = 9003.999899999999  # L11: 9002.25 -> 9003.999899999999
$ws.Cells.Item(11, 13).Value = -421.99999  # M11: -425.5 -> -421.99999
$ws.Cells.Item(11, 14).Value = -9283.999899999999  # N11: -9282.25 -> -9283.999899999999
$ws.Cells.Item(15, 8).Value = 478.83334  # H15: 424.7143 -> 478.83334
$ws.Cells.Item(15, 9).Value = 299.5  # I15: 259.6 -> 299.5
$ws.Cells.Item(15, 11).Value = 898.5  # K15: 778.8000000000001 -> 898.5
$ws.Cells.Item(15, 13).Value = -758.5  # M15: -638.8000000000001 -> -758.5
$ws.Cells.Item(23, 8).Value = 550.6  # H23: 491.5 -> 550.6
$ws.Cells.Item(23, 10).Value = 693  # J23: 444.5 -> 693
$ws.Cells.Item(23, 12).Value = 2079  # L23: 1333.5 -> 2079
$ws.Cells.Item(23, 14).Value = -2549  # N23: -1803.5 -> -2549
$ws.Cells.Item(38, 8).Value = 60.77778  # H38: 61.77778 -> 60.77778
$ws.Cells.Item(38, 9).Value = 70  # I38: 60.166668 -> 70
$ws.Cells.Item(38, 10).Value = 49.25  # J38: 65 -> 49.25
$ws.Cells.Item(38, 11).Value = 210  # K38: 180.500004 -> 210
$ws.Cells.Item(38, 12).Value = 147.75  # L38: 195 -> 147.75
$ws.Cells.Item(38, 13).Value = 137  # M38: 166.499996 -> 137
$ws.Cells.Item(38, 14).Value = -841.75  # N38: -889 -> -841.75
$ws.Cells.Item(60, 8).Value = 3064.3333  # H60: 3085.2222 -> 3064.3333
$ws.Cells.Item(60, 9).Value = 2333  # I60: 2395.6667 -> 2333
$ws.Cells.Item(60, 11).Value = 6999  # K60: 7187.000100000001 -> 6999
$ws.Cells.Item(60, 13).Value = -6748  # M60: -6936.000100000001 -> -6748
$ws.Cells.Item(113, 8).Value = 1350.7858  # H113: 1347.7693 -> 1350.7858
$ws.Cells.Item(113, 10).Value = 1662.3334  # J113: 1696.375 -> 1662.3334
$ws.Cells.Item(113, 12).Value = 4987.0002  # L113: 5089.125 -> 4987.0002
$ws.Cells.Item(113, 14).Value = -9327.0002  # N113: -9429.125 -> -9327.0002
$ws.Cells.Item(131, 8).Value = 2705.4167  # H131: 2540.375 -> 2705.4167
$ws.Cells.Item(131, 9).Value = 2258.25  # I131: 2264.6 -> 2258.25
$ws.Cells.Item(131, 10).Value = 2929  # J131: 3000 -> 2929
$ws.Cells.Item(131, 11).Value = 6774.75  # K131: 6793.799999999999 -> 6774.75
$ws.Cells.Item(131, 12).Value = 8787  # L131: 9000 -> 8787
$ws.Cells.Item(131, 13).Value = -1734.75  # M131: -1753.799999999999 -> -1734.75
$ws.Cells.Item(131, 14).Value = -18867  # N131: -19080 -> -18867
$ws.Cells.Item(132, 8).Value = 1139  # H132: 1356.2858 -> 1139
$ws.Cells.Item(132, 10).Value = 0  # J132: 1899.5 -> 0
$ws.Cells.Item(132, 12).Value = 0  # L132: 17095.5 -> 0
$ws.Cells.Item(132, 14).ClearContents()  # N132: -22155.5 -> (removed)
$ws.Cells.Item(133, 8).Value = 16481.188  # H133: 15599.941 -> 16481.188
$ws.Cells.Item(133, 9).Value = 12171.286  # I133: 10837.375 -> 12171.286
$ws.Cells.Item(133, 11).Value = 36513.858  # K133: 32512.125 -> 36513.858
$ws.Cells.Item(133, 13).Value = -31453.858  # M133: -27452.125 -> -31453.858
$ws.Cells.Item(137, 8).Value = 5503.9165  # H137: 5641.273 -> 5503.9165
$ws.Cells.Item(137, 10).Value = 5719.7  # J137: 5911.5557 -> 5719.7
$ws.Cells.Item(137, 12).Value = 17159.1  # L137: 17734.6671 -> 17159.1
$ws.Cells.Item(137, 14).Value = -27359.1  # N137: -27934.6671 -> -27359.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(19, 8).Value = 0  # H19: 10000 -> 0
$ws.Cells.Item(19, 10).Value = 0  # J19: 10000 -> 0
$ws.Cells.Item(19, 12).Value = 0  # L19: 10000 -> 0
$ws.Cells.Item(19, 14).ClearContents()  # N19: -10576 -> (removed)
$ws.Cells.Item(46, 8).Value = 4444  # H46: 18086 -> 4444
$ws.Cells.Item(46, 9).Value = 4444  # I46: 3722 -> 4444
$ws.Cells.Item(46, 10).Value = 0  # J46: 32450 -> 0
$ws.Cells.Item(46, 11).Value = 4444  # K46: 3722 -> 4444
$ws.Cells.Item(46, 12).Value = 0  # L46: 32450 -> 0
$ws.Cells.Item(46, 13).Value = -4288  # M46: -3566 -> -4288
$ws.Cells.Item(46, 14).ClearContents()  # N46: -32762 -> (removed)
$ws.Cells.Item(102, 8).Value = 2204.9  # H102: 2277.182 -> 2204.9
$ws.Cells.Item(102, 10).Value = 0  # J102: 3000 -> 0
$ws.Cells.Item(102, 12).Value = 0  # L102: 3000 -> 0
$ws.Cells.Item(102, 14).ClearContents()  # N102: -6244 -> (removed)
$ws.Cells.Item(126, 8).Value = 2681.5  # H126: 2779.5715 -> 2681.5
$ws.Cells.Item(126, 10).Value = 3328.3333  # J126: 3995 -> 3328.3333
$ws.Cells.Item(126, 12).Value = 9984.999899999999  # L126: 11985 -> 9984.999899999999
$ws.Cells.Item(126, 14).Value = -14924.9999  # N126: -16925 -> -14924.9999
$ws.Cells.Item(131, 8).Value = 23000  # H131: 0 -> 23000
$ws.Cells.Item(131, 10).Value = 23000  # J131: 0 -> 23000
$ws.Cells.Item(131, 12).Value = 23000  # L131: 0 -> 23000
$ws.Cells.Item(131, 14).Value = -33080  # N131: None -> -33080
$ws.Cells.Item(132, 8).Value = 2170.6667  # H132: 2500 -> 2170.6667
$ws.Cells.Item(132, 9).Value = 2006  # I132: 0 -> 2006
$ws.Cells.Item(132, 11).Value = 6018  # K132: 0 -> 6018
$ws.Cells.Item(132, 13).Value = -3488  # M132: None -> -3488

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3913.818  # H7: 3975.818 -> 3913.818
$ws.Cells.Item(7, 9).Value = 4202.6665  # I7: 4278.4443 -> 4202.6665
$ws.Cells.Item(7, 11).Value = 4202.6665  # K7: 4278.4443 -> 4202.6665
$ws.Cells.Item(7, 13).Value = -4090.6665  # M7: -4166.4443 -> -4090.6665
$ws.Cells.Item(46, 8).Value = 3663.2856  # H46: 3105.2 -> 3663.2856
$ws.Cells.Item(46, 9).Value = 2596.3333  # I46: 1959.6 -> 2596.3333
$ws.Cells.Item(46, 10).Value = 4463.5  # J46: 4250.8 -> 4463.5
$ws.Cells.Item(46, 11).Value = 2596.3333  # K46: 1959.6 -> 2596.3333
$ws.Cells.Item(46, 12).Value = 4463.5  # L46: 4250.8 -> 4463.5
$ws.Cells.Item(46, 13).Value = -2408.3333  # M46: -1771.6 -> -2408.3333
$ws.Cells.Item(46, 14).Value = -4839.5  # N46: -4626.8 -> -4839.5
$ws.Cells.Item(61, 8).Value = 4386.875  # H61: 4497.5 -> 4386.875
$ws.Cells.Item(61, 9).Value = 4315  # I61: 4497.5 -> 4315
$ws.Cells.Item(61, 10).Value = 4602.5  # J61: 0 -> 4602.5
$ws.Cells.Item(61, 11).Value = 4315  # K61: 4497.5 -> 4315
$ws.Cells.Item(61, 12).Value = 4602.5  # L61: 0 -> 4602.5
$ws.Cells.Item(61, 13).Value = -4113  # M61: -4295.5 -> -4113
$ws.Cells.Item(61, 14).Value = -5006.5  # N61: None -> -5006.5
$ws.Cells.Item(82, 8).Value = 4966.3335  # H82: 4950 -> 4966.3335
$ws.Cells.Item(82, 10).Value = 4999  # J82: 0 -> 4999
$ws.Cells.Item(82, 12).Value = 4999  # L82: 0 -> 4999
$ws.Cells.Item(82, 14).Value = -5721  # N82: None -> -5721
$ws.Cells.Item(85, 8).Value = 4966.3335  # H85: 4950 -> 4966.3335
$ws.Cells.Item(85, 10).Value = 4999  # J85: 0 -> 4999
$ws.Cells.Item(85, 12).Value = 4999  # L85: 0 -> 4999
$ws.Cells.Item(85, 14).Value = -7495  # N85: None -> -7495
$ws.Cells.Item(100, 8).Value = 3  # H100: 0 -> 3
$ws.Cells.Item(100, 9).Value = 3  # I100: 0 -> 3
$ws.Cells.Item(100, 11).Value = 3  # K100: 0 -> 3
$ws.Cells.Item(100, 13).Value = 538  # M100: None -> 538
$ws.Cells.Item(113, 8).Value = 4386.875  # H113: 4497.5 -> 4386.875
$ws.Cells.Item(113, 9).Value = 4315  # I113: 4497.5 -> 4315
$ws.Cells.Item(113, 10).Value = 4602.5  # J113: 0 -> 4602.5
$ws.Cells.Item(113, 11).Value = 4315  # K113: 4497.5 -> 4315
$ws.Cells.Item(113, 12).Value = 4602.5  # L113: 0 -> 4602.5
$ws.Cells.Item(113, 13).Value = -2145  # M113: -2327.5 -> -2145
$ws.Cells.Item(113, 14).Value = -8942.5  # N113: None -> -8942.5
$ws.Cells.Item(126, 8).Value = 3913.818  # H126: 3975.818 -> 3913.818
$ws.Cells.Item(126, 9).Value = 4202.6665  # I126: 4278.4443 -> 4202.6665
$ws.Cells.Item(126, 11).Value = 12607.9995  # K126: 12835.3329 -> 12607.9995
$ws.Cells.Item(126, 13).Value = -10137.9995  # M126: -10365.3329 -> -10137.9995
$ws.Cells.Item(140, 8).Value = 99998  # H140: 99997.664 -> 99998
$ws.Cells.Item(140, 10).Value = 99998  # J140: 99997.664 -> 99998
$ws.Cells.Item(140, 12).Value = 99998  # L140: 99997.664 -> 99998
$ws.Cells.Item(140, 14).Value = -110358  # N140: -110357.664 -> -110358

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 31666.666  # H54: 60000 -> 31666.666
$ws.Cells.Item(54, 9).Value = 20000  # I54: 0 -> 20000
$ws.Cells.Item(54, 10).Value = 55000  # J54: 60000 -> 55000
$ws.Cells.Item(54, 11).Value = 20000  # K54: 0 -> 20000
$ws.Cells.Item(54, 12).Value = 55000  # L54: 60000 -> 55000
$ws.Cells.Item(54, 13).Value = -19480  # M54: None -> -19480
$ws.Cells.Item(54, 14).Value = -56040  # N54: -61040 -> -56040
$ws.Cells.Item(81, 8).Value = 1174  # H81: 1521.6666 -> 1174
$ws.Cells.Item(81, 9).Value = 1269.8182  # I81: 1521.6666 -> 1269.8182
$ws.Cells.Item(81, 10).Value = 120  # J81: 0 -> 120
$ws.Cells.Item(81, 11).Value = 2539.6364  # K81: 3043.3332 -> 2539.6364
$ws.Cells.Item(81, 12).Value = 240  # L81: 0 -> 240
$ws.Cells.Item(81, 13).Value = -1478.6364  # M81: -1982.3332 -> -1478.6364
$ws.Cells.Item(81, 14).Value = -2362  # N81: None -> -2362
$ws.Cells.Item(84, 8).Value = 1174  # H84: 1521.6666 -> 1174
$ws.Cells.Item(84, 9).Value = 1269.8182  # I84: 1521.6666 -> 1269.8182
$ws.Cells.Item(84, 10).Value = 120  # J84: 0 -> 120
$ws.Cells.Item(84, 11).Value = 12698.182  # K84: 15216.666 -> 12698.182
$ws.Cells.Item(84, 12).Value = 1200  # L84: 0 -> 1200
$ws.Cells.Item(84, 13).Value = -7394.181999999999  # M84: -9912.666000000001 -> -7394.181999999999
$ws.Cells.Item(84, 14).Value = -11808  # N84: None -> -11808
$ws.Cells.Item(95, 8).Value = 0  # H95: 48650 -> 0
$ws.Cells.Item(95, 10).Value = 0  # J95: 48650 -> 0
$ws.Cells.Item(95, 12).Value = 0  # L95: 48650 -> 0
$ws.Cells.Item(95, 14).ClearContents()  # N95: -54142 -> (removed)
$ws.Cells.Item(96, 8).Value = 900  # H96: 966.3333 -> 900
$ws.Cells.Item(96, 9).Value = 900  # I96: 966.3333 -> 900
$ws.Cells.Item(96, 11).Value = 900  # K96: 966.3333 -> 900
$ws.Cells.Item(96, 13).Value = 473  # M96: 406.6667 -> 473
$ws.Cells.Item(107, 8).Value = 5147.2  # H107: 5497.222 -> 5147.2
$ws.Cells.Item(107, 9).Value = 3927.4285  # I107: 4249.1665 -> 3927.4285
$ws.Cells.Item(107, 11).Value = 11782.2855  # K107: 12747.4995 -> 11782.2855
$ws.Cells.Item(107, 13).Value = -9862.2855  # M107: -10827.4995 -> -9862.2855
$ws.Cells.Item(126, 8).Value = 2074  # H126: 2097.4736 -> 2074
$ws.Cells.Item(126, 9).Value = 2024.0588  # I126: 2050.2942 -> 2024.0588
$ws.Cells.Item(126, 11).Value = 6072.1764  # K126: 6150.882599999999 -> 6072.1764
$ws.Cells.Item(126, 13).Value = -3602.1764  # M126: -3680.882599999999 -> -3602.1764
$ws.Cells.Item(130, 8).Value = 35996.668  # H130: 32747 -> 35996.668
$ws.Cells.Item(130, 10).Value = 35996.668  # J130: 32747 -> 35996.668
$ws.Cells.Item(130, 12).Value = 35996.668  # L130: 32747 -> 35996.668
$ws.Cells.Item(130, 14).Value = -46036.668  # N130: -42787 -> -46036.668
$ws.Cells.Item(132, 8).Value = 4434.6  # H132: 3406.8096 -> 4434.6
$ws.Cells.Item(132, 9).Value = 3393.8572  # I132: 2870.4443 -> 3393.8572
$ws.Cells.Item(132, 10).Value = 19005  # J132: 6625 -> 19005
$ws.Cells.Item(132, 11).Value = 10181.5716  # K132: 8611.332900000001 -> 10181.5716
$ws.Cells.Item(132, 12).Value = 57015  # L132: 19875 -> 57015
$ws.Cells.Item(132, 13).Value = -7651.571599999999  # M132: -6081.332900000001 -> -7651.571599999999
$ws.Cells.Item(132, 14).Value = -62075  # N132: -24935 -> -62075
$ws.Cells.Item(135, 8).Value = 54999  # H135: 112000 -> 54999
$ws.Cells.Item(135, 9).Value = 54999  # I135: 0 -> 54999
$ws.Cells.Item(135, 10).Value = 0  # J135: 112000 -> 0
$ws.Cells.Item(135, 11).Value = 54999  # K135: 0 -> 54999
$ws.Cells.Item(135, 12).Value = 0  # L135: 112000 -> 0
$ws.Cells.Item(135, 13).Value = -49929  # M135: None -> -49929
$ws.Cells.Item(135, 14).ClearContents()  # N135: -122140 -> (removed)
